$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$titles = @(
    "The Itchy & Scratchy & Poochie Show",
    "Homie the Clown",
    "Bart Gets an Elephant",
    "Homer Goes to College",
    "Bart’s Inner Child",
    "Rosebud",
    "Homer’s Enemy",
    "Homer vs. the Eighteenth Amendment",
    "A Milhouse Divided",
    "Marge vs. the Monorail",
    "You Only Move Twice",
    "Homer Badman",
    "Grade School Confidential",
    "Realty Bites",
    "Last Exit to Springfield",
    "The Last Temptation of Homer",
    "The Boy Who Knew Too Much",
    "King-Size Homer",
    "The Cartridge Family ",
    "Deep Space Homer",
    "Homer at the Bat ",
    "Team Homer ",
    "Lisa’s Rival",
    "All’s Fair in Oven War",
    "Homer the Vigilante"
)

for ($i = 0; $i -lt $titles.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $titles[$i]
}
